$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-converted to a number by Excel, so they stay text like the source data.
foreach ($r in @(4,5,7,8,9,10,11,12,14,15,17,18,19,20,21,23,24,25,27,28,29,30,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,47,48,49,50,51)) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.934.79"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.875.88"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.94%  "
$ws.Range("D5").Value = "335.65"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").Value = "0.4765"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.3941"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "47.03"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").Value = "0.07999"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "1.016"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "21.83"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "1.892.39"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "6.047"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "7.188"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "88.46"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "0.06718"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "0.00001049"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "17.04"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").Value = "27.940.81"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "5.499"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "10.99"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "2.339"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "2.107.63"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "158.27"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "19.84"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").Value = "2.101"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "5.456"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "121.37"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Value = "0.9747"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "0.09542"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "3.633"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "5.330"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "1.350"
$ws.Range("E36").Value = "  -6.94%  "
$ws.Range("D37").Value = "0.06076"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "0.02239"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "8.177"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "1.009"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "0.5975"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").Value = "0.1892"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "10.32"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "1.258"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "0.5666"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "12.14"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").Value = "1.926"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "3.334"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "0.06788"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "112.27"
$ws.Range("E51").Value = "  -2.04%  "
